$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.531.67'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '2.217.56'
$ws.Range("E3").Value = '  -2.30%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''311.28'
$ws.Range("E5").Value = '  -2.31%  '

$ws.Range("D6").Value = '''97.01'
$ws.Range("E6").Value = '  -4.92%  '

$ws.Range("D7").Value = '''0.565'
$ws.Range("E7").Value = '  -3.55%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '''0.531'
$ws.Range("E9").Value = '  -6.76%  '

$ws.Range("D10").Value = '''35.39'
$ws.Range("E10").Value = '  -8.48%  '

$ws.Range("D11").Value = '''0.0820'
$ws.Range("E11").Value = '  -2.15%  '

$ws.Range("D12").Value = '''7.31'
$ws.Range("E12").Value = '  -6.75%  '

$ws.Range("E13").Value = '  -3.15%  '

$ws.Range("D14").Value = '2.553.62'
$ws.Range("E14").Value = '  -2.32%  '

$ws.Range("D15").Value = '2.219.62'
$ws.Range("E15").Value = '  -2.23%  '

$ws.Range("D16").Value = '''0.832'
$ws.Range("E16").Value = '  -4.90%  '

$ws.Range("D17").Value = '''13.98'
$ws.Range("E17").Value = '  -3.76%  '

$ws.Range("D18").Value = '43.368.18'
$ws.Range("E18").Value = '  -1.25%  '

$ws.Range("D19").Value = '''12.86'
$ws.Range("E19").Value = '  -9.95%  '

$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -3.07%  '

$ws.Range("D21").Value = '''6.26'
$ws.Range("E21").Value = '  -5.81%  '

$ws.Range("D22").Value = '''65.10'
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").Value = '''233.77'
$ws.Range("E23").Value = '  -1.95%  '

$ws.Range("D24").Value = '''2.95'
$ws.Range("E24").Value = '  -8.04%  '

$ws.Range("D25").Value = '''2.01'
$ws.Range("E25").Value = '  -7.83%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").Value = '''9.93'
$ws.Range("E27").Value = '  -2.78%  '

$ws.Range("D28").Value = '''2.18'
$ws.Range("E28").Value = '  -1.85%  '

$ws.Range("D29").Value = '''35.78'
$ws.Range("E29").Value = '  -8.38%  '

$ws.Range("D30").Value = '''159.81'
$ws.Range("E30").Value = '  -2.12%  '

$ws.Range("D31").Value = '''5.91'
$ws.Range("E31").Value = '  -9.21%  '

$ws.Range("D32").Value = '''19.77'
$ws.Range("E32").Value = '  -3.30%  '

$ws.Range("D33").Value = '''0.0820'
$ws.Range("E33").Value = '  -6.76%  '

$ws.Range("D34").Value = '''2.67'
$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("D35").Value = '''3.08'
$ws.Range("E35").Value = '  -5.19%  '

$ws.Range("D36").Value = '''0.107'
$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("D37").Value = '''1.85'
$ws.Range("E37").Value = '  -9.33%  '

$ws.Range("E38").Value = '  -3.84%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '''3.51'
$ws.Range("E39").Value = '  -9.07%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '''15.09'
$ws.Range("E40").Value = '  -3.99%  '

$ws.Range("D41").Value = '''3.95'
$ws.Range("E41").Value = '  -12.96%  '

$ws.Range("D42").Value = '''0.0305'
$ws.Range("E42").Value = '  -6.27%  '

$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").Value = '1.698.38'
$ws.Range("E44").Value = '  -4.38%  '

$ws.Range("D45").Value = '''82.52'
$ws.Range("E45").Value = '  -2.70%  '

$ws.Range("D46").Value = '''0.192'
$ws.Range("E46").Value = '  -7.02%  '

$ws.Range("D47").Value = '''5.08'
$ws.Range("E47").Value = '  -5.99%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''100.95'
$ws.Range("E48").Value = '  -3.91%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''1.61'
$ws.Range("E49").Value = '  +0.18%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '''55.94'
$ws.Range("E50").Value = '  -6.09%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '''70.04'
$ws.Range("E51").Value = '  -6.21%  '
